$d = $word.ActiveDocument

# Locate the end of the document title ("Lista de necessidades"), i.e. the
# point right before the italic " (Needs list)" subtitle that is about to
# be removed.
$titleRng = $d.Content
$titleRng.Find.Execute("Lista de necessidades", $false, $false, $false, $false, `
                        $false, $true, 1, $false, "", 0)
$titleEnd = $titleRng.End

# Word keeps exactly one "_GoBack" bookmark in the package, marking the
# location of the most recent edit; re-adding it under that reserved name
# automatically drops whatever copy used to sit elsewhere in the document.
# Plant it now, at the spot that is about to become the new "last edit"
# location, before the text below it is deleted.
$goBackRng = $d.Range($titleEnd, $titleEnd)
$d.Bookmarks.Add("_GoBack", $goBackRng)

# Remove the italic " (Needs list)" subtitle that followed the title,
# collapsing the title paragraph back down to just its plain text run.
$d.Content.Find.Execute(" (Needs list)", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)
